$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 14.73304615076469
$ws.Range("C2").Value = 9.384971601072655
$ws.Range("D2").Value = 8.407772783563104
$ws.Range("F2").Value = 39.42189265518708
$ws.Range("G2").Value = 3.701014111491583
$ws.Range("J2").Value = 10.98331999469485
$ws.Range("K2").Value = 10.65598252951231
$ws.Range("L2").Value = 11.31029268606957
$ws.Range("O2").Value = 30.2420889190203
$ws.Range("B3").Value = 14.5374105967684
$ws.Range("C3").Value = 9.386037594419898
$ws.Range("D3").Value = 8.389504438988906
$ws.Range("F3").Value = 39.5148890567641
$ws.Range("G3").Value = 3.702955389642616
$ws.Range("J3").Value = 11.00719526617126
$ws.Range("K3").Value = 10.51175743096998
$ws.Range("L3").Value = 11.30705293847524
$ws.Range("O3").Value = 30.3335505520159
$ws.Range("B4").Value = 14.4185637639663
$ws.Range("C4").Value = 9.387030536253862
$ws.Range("D4").Value = 8.379443085058503
$ws.Range("F4").Value = 39.57963713584184
$ws.Range("G4").Value = 3.704211055586122
$ws.Range("J4").Value = 11.02283462919437
$ws.Range("K4").Value = 10.42381639291468
$ws.Range("L4").Value = 11.30650613382095
$ws.Range("O4").Value = 30.3950376634348
$ws.Range("B5").Value = 14.37050911022455
$ws.Range("C5").Value = 9.387520691711245
$ws.Range("D5").Value = 8.375636533948668
$ws.Range("F5").Value = 39.60794310725539
$ws.Range("G5").Value = 3.704738820850376
$ws.Range("J5").Value = 11.02945467779234
$ws.Range("K5").Value = 10.38817439134601
$ws.Range("L5").Value = 11.30664741545261
$ws.Range("O5").Value = 30.42143298267713
$ws.Range("B6").Value = 14.36255403228694
$ws.Range("C6").Value = 9.387607259529874
$ws.Range("D6").Value = 8.375022270631563
$ws.Range("F6").Value = 39.61275921308401
$ws.Range("G6").Value = 3.704827428001353
$ws.Range("J6").Value = 11.03056885727623
$ws.Range("K6").Value = 10.3822689992229
$ws.Range("L6").Value = 11.30669290507366
$ws.Range("O6").Value = 30.42589672873944
$ws.Range("B7").Value = 14.41791408444694
$ws.Range("C7").Value = 9.387036799815197
$ws.Range("D7").Value = 8.379390556289758
$ws.Range("F7").Value = 39.58001110716985
$ws.Range("G7").Value = 3.704218108075198
$ws.Range("J7").Value = 11.02292290921951
$ws.Range("K7").Value = 10.42333487092989
$ws.Range("L7").Value = 11.30650656316364
$ws.Range("O7").Value = 30.39538822083781
$ws.Range("B8").Value = 14.66536025250827
$ws.Range("C8").Value = 9.385269171438001
$ws.Range("D8").Value = 8.401235882147489
$ws.Range("F8").Value = 39.45236904519739
$ws.Range("G8").Value = 3.701670270625943
$ws.Range("J8").Value = 10.99134911523801
$ws.Range("K8").Value = 10.60615055795536
$ws.Range("L8").Value = 11.30887712238935
$ws.Range("O8").Value = 30.2725180445248
$ws.Range("B9").Value = 15.1582419223926
$ws.Range("C9").Value = 9.384469815530132
$ws.Range("D9").Value = 8.45310408146117
$ws.Range("F9").Value = 39.26285063730114
$ws.Range("G9").Value = 3.697177225467671
$ws.Range("J9").Value = 10.93718562773574
$ws.Range("K9").Value = 10.96775105111653
$ws.Range("L9").Value = 11.32490537462068
$ws.Range("O9").Value = 30.07391110549016
$ws.Range("B10").Value = 15.52172718927837
$ws.Range("C10").Value = 9.385484558316005
$ws.Range("D10").Value = 8.496519227679455
$ws.Range("F10").Value = 39.16079914003178
$ws.Range("G10").Value = 3.694179791756575
$ws.Range("J10").Value = 10.90208752357299
$ws.Range("K10").Value = 11.23296008093223
$ws.Range("L10").Value = 11.3435255410902
$ws.Range("O10").Value = 29.95388056379972
$ws.Range("B11").Value = 15.68669760747831
$ws.Range("C11").Value = 9.38628901737577
$ws.Range("D11").Value = 8.517376449899736
$ws.Range("F11").Value = 39.12247025462251
$ws.Range("G11").Value = 3.692881428128455
$ws.Range("J11").Value = 10.88713373065079
$ws.Range("K11").Value = 11.35302551631791
$ws.Range("L11").Value = 11.35346016371427
$ws.Range("O11").Value = 29.90491123295025
$ws.Range("B12").Value = 15.74905721766733
$ws.Range("C12").Value = 9.38664251079326
$ws.Range("D12").Value = 8.525429586304991
$ws.Range("F12").Value = 39.10912137723955
$ws.Range("G12").Value = 3.692399093602158
$ws.Range("J12").Value = 10.881616230817
$ws.Range("K12").Value = 11.39836858448287
$ws.Range("L12").Value = 11.35743056726014
$ws.Range("O12").Value = 29.8871790217765
$ws.Range("B13").Value = 15.73563278103894
$ws.Range("C13").Value = 9.38656421350154
$ws.Range("D13").Value = 8.523688377626188
$ws.Range("F13").Value = 39.1119444448266
$ws.Range("G13").Value = 3.692502558849941
$ws.Range("J13").Value = 10.88279807299463
$ws.Range("K13").Value = 11.38860923855105
$ws.Range("L13").Value = 11.35656623816477
$ws.Range("O13").Value = 29.89096186195961
$ws.Range("B14").Value = 15.69183051380932
$ws.Range("C14").Value = 9.386317122738189
$ws.Range("D14").Value = 8.518035905867247
$ws.Range("F14").Value = 39.12134867045647
$ws.Range("G14").Value = 3.692841559498823
$ws.Range("J14").Value = 10.8866768954497
$ws.Range("K14").Value = 11.35675860842258
$ws.Range("L14").Value = 11.35378264791958
$ws.Range("O14").Value = 29.90343612378866
$ws.Range("B15").Value = 15.6649842464752
$ws.Range("C15").Value = 9.386172122770597
$ws.Range("D15").Value = 8.514593653602965
$ws.Range("F15").Value = 39.12726084077345
$ws.Range("G15").Value = 3.693050420561322
$ws.Range("J15").Value = 10.8890716801213
$ws.Range("K15").Value = 11.33723198655801
$ws.Range("L15").Value = 11.3521046886929
$ws.Range("O15").Value = 29.91118267846388
$ws.Range("B16").Value = 15.51093346846181
$ws.Range("C16").Value = 9.385438846509309
$ws.Range("D16").Value = 8.495178081125554
$ws.Range("F16").Value = 39.16346704034272
$ws.Range("G16").Value = 3.694265950620415
$ws.Range("J16").Value = 10.90308512733934
$ws.Range("K16").Value = 11.22509843330523
$ws.Range("L16").Value = 11.34290557092703
$ws.Range("O16").Value = 29.95719429715954
$ws.Range("B17").Value = 15.41629018126615
$ws.Range("C17").Value = 9.385076504905829
$ws.Range("D17").Value = 8.483547904637094
$ws.Range("F17").Value = 39.18775276910865
$ws.Range("G17").Value = 3.695028301434873
$ws.Range("J17").Value = 10.91194095076637
$ws.Range("K17").Value = 11.15613102494093
$ws.Range("L17").Value = 11.33763572440703
$ws.Range("O17").Value = 29.98686472120495
$ws.Range("B18").Value = 15.36182129716611
$ws.Range("C18").Value = 9.384900392347564
$ws.Range("D18").Value = 8.476963041863513
$ws.Range("F18").Value = 39.20248314093994
$ws.Range("G18").Value = 3.695472923598869
$ws.Range("J18").Value = 10.91712991334691
$ws.Range("K18").Value = 11.11641055945433
$ws.Range("L18").Value = 11.33474260530549
$ws.Range("O18").Value = 30.00446050706187
$ws.Range("B19").Value = 15.34337521759615
$ws.Range("C19").Value = 9.384846324007972
$ws.Range("D19").Value = 8.474751597539763
$ws.Range("F19").Value = 39.20760139388511
$ws.Range("G19").Value = 3.695624520656756
$ws.Range("J19").Value = 10.91890319234838
$ws.Range("K19").Value = 11.102954152704
$ws.Range("L19").Value = 11.33378680332884
$ws.Range("O19").Value = 30.01050914743382
$ws.Range("B20").Value = 15.42636890464821
$ws.Range("C20").Value = 9.385111737372064
$ws.Range("D20").Value = 8.484775171220438
$ws.Range("F20").Value = 39.18508865053511
$ws.Range("G20").Value = 3.694946512947098
$ws.Range("J20").Value = 10.91098837030188
$ws.Range("K20").Value = 11.16347842162326
$ws.Range("L20").Value = 11.33818244585868
$ws.Range("O20").Value = 29.98365137595829
$ws.Range("B21").Value = 15.70469975781747
$ws.Range("C21").Value = 9.38638837668187
$ws.Range("D21").Value = 8.519692003360227
$ws.Range("F21").Value = 39.11855478221261
$ws.Range("G21").Value = 3.692741734005355
$ws.Range("J21").Value = 10.88553365454
$ws.Range("K21").Value = 11.36611756239364
$ws.Range("L21").Value = 11.35459461851148
$ws.Range("O21").Value = 29.89975009868594
$ws.Range("B22").Value = 15.88593250949117
$ws.Range("C22").Value = 9.387507336251174
$ws.Range("D22").Value = 8.543413280335201
$ws.Range("F22").Value = 39.08186453454498
$ws.Range("G22").Value = 3.691355130935524
$ws.Range("J22").Value = 10.86974353688519
$ws.Range("K22").Value = 11.49781886674213
$ws.Range("L22").Value = 11.3665344069359
$ws.Range("O22").Value = 29.84964587366019
$ws.Range("B23").Value = 15.78928489878637
$ws.Range("C23").Value = 9.386884230672132
$ws.Range("D23").Value = 8.530671798990991
$ws.Range("F23").Value = 39.10082483056429
$ws.Range("G23").Value = 3.692090229324637
$ws.Range("J23").Value = 10.87809374770056
$ws.Range("K23").Value = 11.42760747105744
$ws.Range("L23").Value = 11.36005163476707
$ws.Range("O23").Value = 29.87595420322515
$ws.Range("B24").Value = 15.42181248765331
$ws.Range("C24").Value = 9.385095708435717
$ws.Range("D24").Value = 8.484220007587357
$ws.Range("F24").Value = 39.18629070579399
$ws.Range("G24").Value = 3.694983469780692
$ws.Range("J24").Value = 10.91141872778278
$ws.Range("K24").Value = 11.1601568791458
$ws.Range("L24").Value = 11.3379348474826
$ws.Range("O24").Value = 29.98510245401103
$ws.Range("B25").Value = 15.02442855324919
$ws.Range("C25").Value = 9.384403002746392
$ws.Range("D25").Value = 8.438125235978216
$ws.Range("F25").Value = 39.30759806970661
$ws.Range("G25").Value = 3.698339167523001
$ws.Range("J25").Value = 10.95101147347354
$ws.Range("K25").Value = 11.38860923855105
$ws.Range("L25").Value = 11.35656623816477
$ws.Range("O25").Value = 30.12309879538133
